# Update countries & provincias Spain
#
# 1) Reorders four pairs of countries in the "Pais" sheet so that the
#    row order matches the new canonical list (Venezuela before Estado
#    de Palestina, Santa Lucia before Belice, Namibia before San
#    Vicente y las Granadinas, Burundi before San Cristobal y Nieves),
#    and
# 2) Refreshes the case-count figures for the USA row (row 4) and for
#    the rows whose country labels moved (Venezuela gets new figures;
#    the swapped pairs keep their original figures, just relocated).
#
# Each affected row is rewritten in full (columns A:H) with its final
# target values so the result is correct regardless of the workbook's
# starting state. (Note: this engine's PowerShell function parameter
# binding only works positionally, not via -Name value syntax.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CountryRow($Row, $Pais, $CasosTotales, $NuevosCasos, $CasosActivos, $Recuperados, $CasosCriticos, $MuertesHoy, $Muertes) {
    $ws.Range("A$Row").Value = $Pais
    $ws.Range("B$Row").Value = $CasosTotales
    $ws.Range("C$Row").Value = $NuevosCasos
    $ws.Range("D$Row").Value = $CasosActivos
    $ws.Range("E$Row").Value = $Recuperados
    $ws.Range("F$Row").Value = $CasosCriticos
    $ws.Range("G$Row").Value = $MuertesHoy
    $ws.Range("H$Row").Value = $Muertes
}

# Row 4: Estados Unidos - refreshed totals
Set-CountryRow 4 "Estados Unidos" 1188029 27255 178263 941179 16139 1143 68587

# Rows 124/125: Venezuela now precedes Estado de Palestina; Venezuela gets new figures
Set-CountryRow 124 "Venezuela" 357 12 148 199 2 0 10
Set-CountryRow 125 "Estado de Palestina" 353 0 77 274 0 0 2

# Rows 188/189: Santa Lucia now precedes Belice (rows swap places, figures unchanged)
Set-CountryRow 188 "Santa Lucia" 18 1 15 3 0 0 0
Set-CountryRow 189 "Belice" 18 0 13 3 1 0 2

# Rows 194/195: Namibia now precedes San Vicente y las Granadinas (identical figures, labels swap)
Set-CountryRow 194 "Namibia" 16 0 8 8 0 0 0
Set-CountryRow 195 "San Vicente y las Granadinas" 16 0 8 8 0 0 0

# Rows 198/199: Burundi now precedes San Cristobal y Nieves (rows swap places, figures unchanged)
Set-CountryRow 198 "Burundi" 15 0 7 7 0 0 1
Set-CountryRow 199 "San Cristobal y Nieves" 15 0 8 7 0 0 0

Write-Host "Countries reordered and figures updated."
